$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, shifting existing rows 250:264 down to 251:265.
$ws.Rows("250:250").Insert()

# Populate the newly inserted row 250 with the latest weekly record.
$ws.Range("A250").Value = 11
$ws.Range("B250").Value = "Vega Monumental Concepción"
$ws.Range("C250").Value = "Bíobío"
$ws.Range("D250").Value = 45013
$ws.Range("D250").NumberFormat = $ws.Range("D251").NumberFormat
$ws.Range("E250").Value = 8
$ws.Range("F250").Value = 100112003
$ws.Range("G250").Value = "Ajo"
$ws.Range("H250").Value = "Chino"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 270
$ws.Range("K250").Value = 14000
$ws.Range("L250").Value = 15000
$ws.Range("M250").Value = 14556
$ws.Range("N250").Value = "$/caja 10 kilos"
$ws.Range("O250").Value = "China"
$ws.Range("P250").Value = 1456
$ws.Range("Q250").Value = 10
$ws.Range("R250").Value = "Hortaliza"
